# Apply updates to the pharmacy stock report workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 corresponds to "VOLTAREN 75MG/3ML 3 AMP."
$ws.Range("H21").Value = "3:2"
$ws.Range("L21").Value = 51
$ws.Range("N21").Value = "1:0"

# Row 25 corresponds to "سرنجات 3 سم"
$ws.Range("H25").Value = "-1:0"
$ws.Range("L25").Value = 12
$ws.Range("N25").Value = "6:0"

# Update the cached total in K35 (sum of price column L4:L34) to reflect the
# increases above (34 -> 51 and 10 -> 12, i.e. +19 total).
$ws.Range("K35").Value = 1464.58
